$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B, rows 2:25
$Bvals = New-Object 'object[,]' 24,1
$Bvals[0,0] = 0.1978249331691586
$Bvals[1,0] = 0.1737139522058726
$Bvals[2,0] = 0.1588884327514535
$Bvals[3,0] = 0.1528420462135358
$Bvals[4,0] = 0.1518377686176962
$Bvals[5,0] = 0.1588069081331582
$Bvals[6,0] = 0.1895161665353839
$Bvals[7,0] = 0.2495496360358516
$Bvals[8,0] = 0.2935213868330777
$Bvals[9,0] = 0.313491849838698
$Bvals[10,0] = 0.3210490262205496
$Bvals[11,0] = 0.3194216929870777
$Bvals[12,0] = 0.3141136906323538
$Bvals[13,0] = 0.3108616924027388
$Bvals[14,0] = 0.2922155706775698
$Bvals[15,0] = 0.2807680760341498
$Bvals[16,0] = 0.2741807518741837
$Bvals[17,0] = 0.2719498931844839
$Bvals[18,0] = 0.2819869986083745
$Bvals[19,0] = 0.3156729250671617
$Bvals[20,0] = 0.3376580430853551
$Bvals[21,0] = 0.325927145326034
$Bvals[22,0] = 0.2814359423828137
$Bvals[23,0] = 0.2333311788004266
$ws.Range("B2:B25").Value = $Bvals

# Column C, rows 2:25
$Cvals = New-Object 'object[,]' 24,1
$Cvals[0,0] = 0.03345123140478279
$Cvals[1,0] = 0.03034255889001969
$Cvals[2,0] = 0.02841830620812402
$Cvals[3,0] = 0.02763030913501296
$Cvals[4,0] = 0.02749923193683657
$Cvals[5,0] = 0.02840769451668734
$Cvals[6,0] = 0.0323826122230173
$Cvals[7,0] = 0.04005230873728749
$Cvals[8,0] = 0.04560876030262762
$Cvals[9,0] = 0.04811905015962736
$Cvals[10,0] = 0.04906708500452339
$Cvals[11,0] = 0.04886302343879834
$Cvals[12,0] = 0.04819709715634701
$Cvals[13,0] = 0.04778886324230314
$Cvals[14,0] = 0.04544435234842581
$Cvals[15,0] = 0.04400157993099185
$Cvals[16,0] = 0.04317010425130263
$Cvals[17,0] = 0.04288830254674281
$Cvals[18,0] = 0.04415533476999656
$Cvals[19,0] = 0.04839276573470386
$Cvals[20,0] = 0.05114723710001101
$Cvals[21,0] = 0.049678510294072
$Cvals[22,0] = 0.04408582845312026
$Cvals[23,0] = 0.03799107353489717
$ws.Range("C2:C25").Value = $Cvals

# Column D, rows 2:25
$Dvals = New-Object 'object[,]' 24,1
$Dvals[0,0] = 0.02829542261005713
$Dvals[1,0] = 0.02592591746183359
$Dvals[2,0] = 0.02445805973165704
$Dvals[3,0] = 0.02385667140514158
$Dvals[4,0] = 0.02375661768862614
$Dvals[5,0] = 0.02444996220370399
$Dvals[6,0] = 0.02748113146429887
$Dvals[7,0] = 0.03332091536897508
$Dvals[8,0] = 0.03754627809394862
$Dvals[9,0] = 0.03945405813947644
$Dvals[10,0] = 0.04017438537842821
$Dvals[11,0] = 0.04001934444307409
$Dvals[12,0] = 0.03951336238236536
$Dvals[13,0] = 0.0392031580036587
$Dvals[14,0] = 0.03742130758463702
$Dvals[15,0] = 0.03632449199729137
$Dvals[16,0] = 0.03569228375172884
$Dvals[17,0] = 0.03547799852766076
$Dvals[18,0] = 0.03644138980865108
$Dvals[19,0] = 0.03966203912338528
$Dvals[20,0] = 0.04175461453917251
$Dvals[21,0] = 0.04063890779431034
$Dvals[22,0] = 0.03638854539546799
$Dvals[23,0] = 0.03175242535902356
$ws.Range("D2:D25").Value = $Dvals

# Column F, rows 2:25
$Fvals = New-Object 'object[,]' 24,1
$Fvals[0,0] = 0.6317922821469537
$Fvals[1,0] = 0.6297349220110533
$Fvals[2,0] = 0.6288431017936205
$Fvals[3,0] = 0.6285730735984885
$Fvals[4,0] = 0.6285338765906232
$Fvals[5,0] = 0.6288390819419831
$Fvals[6,0] = 0.6310058176044535
$Fvals[7,0] = 0.6382028155314501
$Fvals[8,0] = 0.6452908855658706
$Fvals[9,0] = 0.6489071961430071
$Fvals[10,0] = 0.6503329894089163
$Fvals[11,0] = 0.6500234117293502
$Fvals[12,0] = 0.6490233670685441
$Fvals[13,0] = 0.6484181530727113
$Fvals[14,0] = 0.6450624390062458
$Fvals[15,0] = 0.6431042081978262
$Fvals[16,0] = 0.6420147731246999
$Fvals[17,0] = 0.6416522440298067
$Fvals[18,0] = 0.643308847363123
$Fvals[19,0] = 0.6493155742781695
$Fvals[20,0] = 0.6535699226587184
$Fvals[21,0] = 0.6512692232640731
$Fvals[22,0] = 0.6432162166887849
$Fvals[23,0] = 0.6359398963555236
$ws.Range("F2:F25").Value = $Fvals

# Column G, rows 2:25
$Gvals = New-Object 'object[,]' 24,1
$Gvals[0,0] = 0.002418114716083965
$Gvals[1,0] = 0.002420221227828925
$Gvals[2,0] = 0.002421583355402657
$Gvals[3,0] = 0.002422155768602624
$Gvals[4,0] = 0.002422251865909277
$Gvals[5,0] = 0.002421591005034242
$Gvals[6,0] = 0.00241882681127712
$Gvals[7,0] = 0.002413949008642973
$Gvals[8,0] = 0.002410692671786063
$Gvals[9,0] = 0.002409281620339467
$Gvals[10,0] = 0.002408757339986988
$Gvals[11,0] = 0.002408869806655479
$Gvals[12,0] = 0.002409238286213458
$Gvals[13,0] = 0.002409465298310034
$Gvals[14,0] = 0.002410786297129384
$Gvals[15,0] = 0.002411614649203609
$Gvals[16,0] = 0.002412097712934483
$Gvals[17,0] = 0.00241226240813556
$Gvals[18,0] = 0.002411525785280364
$Gvals[19,0] = 0.002409129782351326
$Gvals[20,0] = 0.00240762243675386
$Gvals[21,0] = 0.002408421591990981
$Gvals[22,0] = 0.002411565939420855
$Gvals[23,0] = 0.002415210839241654
$ws.Range("G2:G25").Value = $Gvals

# Column I, rows 2:25
$Ivals = New-Object 'object[,]' 24,1
$Ivals[0,0] = 0.4910599711396344
$Ivals[1,0] = 0.494195492585348
$Ivals[2,0] = 0.4963726320953299
$Ivals[3,0] = 0.4973231528854427
$Ivals[4,0] = 0.4974848096360134
$Ivals[5,0] = 0.4963851948151969
$Ivals[6,0] = 0.4920887964555902
$Ivals[7,0] = 0.4856639078779601
$Ivals[8,0] = 0.4821654118750018
$Ivals[9,0] = 0.4808397314588717
$Ivals[10,0] = 0.4803759964467957
$Ivals[11,0] = 0.4804741671563697
$Ivals[12,0] = 0.4808008124003074
$Ivals[13,0] = 0.4810058775814881
$Ivals[14,0] = 0.482257401037387
$Ivals[15,0] = 0.4830932834778103
$Ivals[16,0] = 0.4835990766044738
$Ivals[17,0] = 0.4837746244362613
$Ivals[18,0] = 0.4830017129860025
$Ivals[19,0] = 0.4807038297035966
$Ivals[20,0] = 0.47942512843602
$Ivals[21,0] = 0.4800871654034822
$Ivals[22,0] = 0.4830430334133844
$Ivals[23,0] = 0.4871876139626927
$ws.Range("I2:I25").Value = $Ivals

# Column K, rows 2:25
$Kvals = New-Object 'object[,]' 24,1
$Kvals[0,0] = 0.1989937502063128
$Kvals[1,0] = 0.1736316488547374
$Kvals[2,0] = 0.158011040340881
$Kvals[3,0] = 0.1516338889658329
$Kvals[4,0] = 0.1505742799610488
$Kvals[5,0] = 0.157925082277643
$Kvals[6,0] = 0.1902591815677681
$Kvals[7,0] = 0.2532654661944775
$Kvals[8,0] = 0.2992911899580974
$Kvals[9,0] = 0.3201677967536227
$Kvals[10,0] = 0.3280640443915672
$Kvals[11,0] = 0.3263638670504463
$Kvals[12,0] = 0.3208176146727624
$Kvals[13,0] = 0.3174191517069289
$Kvals[14,0] = 0.2979255866727328
$Kvals[15,0] = 0.2859509799764908
$Kvals[16,0] = 0.2790578193080933
$Kvals[17,0] = 0.2767229537657272
$Kvals[18,0] = 0.2872262894287587
$Kvals[19,0] = 0.3224469393219636
$Kvals[20,0] = 0.3454113857322056
$Kvals[21,0] = 0.333159973418816
$Kvals[22,0] = 0.2866497493938027
$Kvals[23,0] = 0.2362656388399671
$ws.Range("K2:K25").Value = $Kvals

# Column M, rows 2:25
$Mvals = New-Object 'object[,]' 24,1
$Mvals[0,0] = 0.8875589759514355
$Mvals[1,0] = 0.790241033318793
$Mvals[2,0] = 0.7310048131484592
$Mvals[3,0] = 0.7069898664091596
$Mvals[4,0] = 0.7030095419686688
$Mvals[5,0] = 0.7306804434995371
$Mvals[6,0] = 0.853892812363668
$Mvals[7,0] = 1.099894404191076
$Mvals[8,0] = 1.283709821964734
$Mvals[9,0] = 1.368091651637599
$Mvals[10,0] = 1.40016188746435
$Mvals[11,0] = 1.393249687678505
$Mvals[12,0] = 1.370727710843127
$Mvals[13,0] = 1.35694775137361
$Mvals[14,0] = 1.278211216853464
$Mvals[15,0] = 1.23010925987829
$Mvals[16,0] = 1.202513706047213
$Mvals[17,0] = 1.193182384261618
$Mvals[18,0] = 1.235222354138529
$Mvals[19,0] = 1.377339739123727
$Mvals[20,0] = 1.47090502708194
$Mvals[21,0] = 1.420902616640547
$Mvals[22,0] = 1.232910541277477
$Mvals[23,0] = 1.032833813049137
$ws.Range("M2:M25").Value = $Mvals

# Column N, rows 2:25
$Nvals = New-Object 'object[,]' 24,1
$Nvals[0,0] = 1.368790931256127
$Nvals[1,0] = 1.385109425173036
$Nvals[2,0] = 1.395621479854954
$Nvals[3,0] = 1.400028982284982
$Nvals[4,0] = 1.400768317136775
$Nvals[5,0] = 1.395680420002756
$Nvals[6,0] = 1.374315301770745
$Nvals[7,0] = 1.33632813621848
$Nvals[8,0] = 1.310804957828459
$Nvals[9,0] = 1.299712453130891
$Nvals[10,0] = 1.295586602515842
$Nvals[11,0] = 1.2964718564241
$Nvals[12,0] = 1.299371520229004
$Nvals[13,0] = 1.301157372472353
$Nvals[14,0] = 1.311540316443624
$Nvals[15,0] = 1.318042713783958
$Nvals[16,0] = 1.321831480077567
$Nvals[17,0] = 1.323122661624245
$Nvals[18,0] = 1.317345475113054
$Nvals[19,0] = 1.298517791868578
$Nvals[20,0] = 1.286647937999993
$Nvals[21,0] = 1.292943242642462
$Nvals[22,0] = 1.317660539581533
$Nvals[23,0] = 1.346185859688617
$ws.Range("N2:N25").Value = $Nvals

# Column O, rows 2:25
$Ovals = New-Object 'object[,]' 24,1
$Ovals[0,0] = 2.142349206341578
$Ovals[1,0] = 2.147525867525161
$Ovals[2,0] = 2.151877565559445
$Ovals[3,0] = 2.153945840081263
$Ovals[4,0] = 2.154307085492846
$Ovals[5,0] = 2.151904265033281
$Ovals[6,0] = 2.143890560608128
$Ovals[7,0] = 2.13749095671335
$Ovals[8,0] = 2.138480157487777
$Ovals[9,0] = 2.140168811174817
$Ovals[10,0] = 2.140986564287346
$Ovals[11,0] = 2.14080251385775
$Ovals[12,0] = 2.14023251385396
$Ovals[13,0] = 2.13990659722009
$Ovals[14,0] = 2.138394740238141
$Ovals[15,0] = 2.137784642055294
$Ovals[16,0] = 2.137550307575395
$Ovals[17,0] = 2.137490982892672
$Ovals[18,0] = 2.13783752140111
$Ovals[19,0] = 2.140395096315672
$Ovals[20,0] = 2.143105953276319
$Ovals[21,0] = 2.141563959032226
$Ovals[22,0] = 2.137813252042747
$Ovals[23,0] = 2.138223709925541
$ws.Range("O2:O25").Value = $Ovals
